# Updates cryptos list (Coin / Link / Price / Volume(1h)) to the latest
# scraped snapshot. Cells D2:E51 / B27:C28 / B38:C39 change value; some
# Price cells are plain numeric-looking text (e.g. "5.59") that Excel
# would otherwise auto-convert to a Double, so those are forced back to
# Text ("@") number format before the value is written, exactly as a
# human re-typing the figure into a text-formatted column would do.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Ref = 'D2'; Value = '63.008.50'; AsText = $false },
    @{ Ref = 'E2'; Value = '  -0.45%  '; AsText = $false },
    @{ Ref = 'D3'; Value = '2.562.75'; AsText = $false },
    @{ Ref = 'E3'; Value = '  +4.12%  '; AsText = $false },
    @{ Ref = 'E4'; Value = '  +0.07%  '; AsText = $false },
    @{ Ref = 'D5'; Value = '568.93'; AsText = $true },
    @{ Ref = 'E5'; Value = '  +0.04%  '; AsText = $false },
    @{ Ref = 'D6'; Value = '147.72'; AsText = $true },
    @{ Ref = 'E6'; Value = '  +2.88%  '; AsText = $false },
    @{ Ref = 'E7'; Value = '  +0.07%  '; AsText = $false },
    @{ Ref = 'D8'; Value = '0.580'; AsText = $true },
    @{ Ref = 'E8'; Value = '  -1.63%  '; AsText = $false },
    @{ Ref = 'D9'; Value = '2.562.23'; AsText = $false },
    @{ Ref = 'E9'; Value = '  +4.15%  '; AsText = $false },
    @{ Ref = 'E10'; Value = '  -1.09%  '; AsText = $false },
    @{ Ref = 'D11'; Value = '5.59'; AsText = $true },
    @{ Ref = 'E11'; Value = '  -3.08%  '; AsText = $false },
    @{ Ref = 'E12'; Value = '  +0.13%  '; AsText = $false },
    @{ Ref = 'E13'; Value = '  +0.15%  '; AsText = $false },
    @{ Ref = 'D14'; Value = '27.16'; AsText = $true },
    @{ Ref = 'E14'; Value = '  +2.69%  '; AsText = $false },
    @{ Ref = 'D15'; Value = '3.022.21'; AsText = $false },
    @{ Ref = 'E15'; Value = '  +4.47%  '; AsText = $false },
    @{ Ref = 'D16'; Value = '63.035.16'; AsText = $false },
    @{ Ref = 'E16'; Value = '  -0.19%  '; AsText = $false },
    @{ Ref = 'D17'; Value = '0.0000142'; AsText = $true },
    @{ Ref = 'E17'; Value = '  -1.56%  '; AsText = $false },
    @{ Ref = 'D18'; Value = '2.519.83'; AsText = $false },
    @{ Ref = 'E18'; Value = '  +2.83%  '; AsText = $false },
    @{ Ref = 'D19'; Value = '11.51'; AsText = $true },
    @{ Ref = 'E19'; Value = '  +1.53%  '; AsText = $false },
    @{ Ref = 'D20'; Value = '334.89'; AsText = $true },
    @{ Ref = 'E20'; Value = '  -2.70%  '; AsText = $false },
    @{ Ref = 'D21'; Value = '4.29'; AsText = $true },
    @{ Ref = 'E21'; Value = '  -0.28%  '; AsText = $false },
    @{ Ref = 'D22'; Value = '6.78'; AsText = $true },
    @{ Ref = 'E22'; Value = '  -1.31%  '; AsText = $false },
    @{ Ref = 'E23'; Value = '  -0.08%  '; AsText = $false },
    @{ Ref = 'D24'; Value = '65.03'; AsText = $true },
    @{ Ref = 'E24'; Value = '  -0.79%  '; AsText = $false },
    @{ Ref = 'E25'; Value = '  -3.28%  '; AsText = $false },
    @{ Ref = 'D26'; Value = '1.60'; AsText = $true },
    @{ Ref = 'E26'; Value = '  +4.42%  '; AsText = $false },
    @{ Ref = 'B27'; Value = 'Binance-PegBSC-USD'; AsText = $false },
    @{ Ref = 'C27'; Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; AsText = $false },
    @{ Ref = 'D27'; Value = '0.999'; AsText = $true },
    @{ Ref = 'E27'; Value = '  -0.10%  '; AsText = $false },
    @{ Ref = 'B28'; Value = 'SuiNetwork'; AsText = $false },
    @{ Ref = 'C28'; Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'; AsText = $false },
    @{ Ref = 'D28'; Value = '1.48'; AsText = $true },
    @{ Ref = 'E28'; Value = '  +11.14%  '; AsText = $false },
    @{ Ref = 'D29'; Value = '8.40'; AsText = $true },
    @{ Ref = 'E29'; Value = '  +1.78%  '; AsText = $false },
    @{ Ref = 'D30'; Value = '7.24'; AsText = $true },
    @{ Ref = 'E30'; Value = '  +5.00%  '; AsText = $false },
    @{ Ref = 'D31'; Value = '0.0₃0815'; AsText = $false },
    @{ Ref = 'E31'; Value = '  -0.43%  '; AsText = $false },
    @{ Ref = 'E32'; Value = '  +1.71%  '; AsText = $false },
    @{ Ref = 'D33'; Value = '177.03'; AsText = $true },
    @{ Ref = 'E33'; Value = '  +1.36%  '; AsText = $false },
    @{ Ref = 'D34'; Value = '1.58'; AsText = $true },
    @{ Ref = 'E34'; Value = '  +5.46%  '; AsText = $false },
    @{ Ref = 'D35'; Value = '412.53'; AsText = $true },
    @{ Ref = 'E35'; Value = '  +11.72%  '; AsText = $false },
    @{ Ref = 'D36'; Value = '0.398'; AsText = $true },
    @{ Ref = 'E36'; Value = '  -0.61%  '; AsText = $false },
    @{ Ref = 'D37'; Value = '18.84'; AsText = $true },
    @{ Ref = 'E37'; Value = '  -0.12%  '; AsText = $false },
    @{ Ref = 'B38'; Value = 'USDe'; AsText = $false },
    @{ Ref = 'C38'; Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'; AsText = $false },
    @{ Ref = 'D38'; Value = '0.999'; AsText = $true },
    @{ Ref = 'E38'; Value = '  +0.00%  '; AsText = $false },
    @{ Ref = 'B39'; Value = 'NEARProtocol'; AsText = $false },
    @{ Ref = 'C39'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; AsText = $false },
    @{ Ref = 'D39'; Value = '4.37'; AsText = $true },
    @{ Ref = 'E39'; Value = '  -2.94%  '; AsText = $false },
    @{ Ref = 'E40'; Value = '  +1.69%  '; AsText = $false },
    @{ Ref = 'D41'; Value = '1.00'; AsText = $true },
    @{ Ref = 'E41'; Value = '  +0.00%  '; AsText = $false },
    @{ Ref = 'D42'; Value = '39.19'; AsText = $true },
    @{ Ref = 'E42'; Value = '  -2.47%  '; AsText = $false },
    @{ Ref = 'D43'; Value = '152.04'; AsText = $true },
    @{ Ref = 'E43'; Value = '  +0.71%  '; AsText = $false },
    @{ Ref = 'D44'; Value = '3.74'; AsText = $true },
    @{ Ref = 'E44'; Value = '  +0.09%  '; AsText = $false },
    @{ Ref = 'D45'; Value = '20.67'; AsText = $true },
    @{ Ref = 'E45'; Value = '  -0.61%  '; AsText = $false },
    @{ Ref = 'D46'; Value = '0.606'; AsText = $true },
    @{ Ref = 'E46'; Value = '  +1.31%  '; AsText = $false },
    @{ Ref = 'D47'; Value = '0.0960'; AsText = $true },
    @{ Ref = 'E47'; Value = '  -0.41%  '; AsText = $false },
    @{ Ref = 'D48'; Value = '0.0520'; AsText = $true },
    @{ Ref = 'E48'; Value = '  -0.79%  '; AsText = $false },
    @{ Ref = 'D49'; Value = '0.0237'; AsText = $true },
    @{ Ref = 'E49'; Value = '  +4.93%  '; AsText = $false },
    @{ Ref = 'D50'; Value = '18.43'; AsText = $true },
    @{ Ref = 'E50'; Value = '  +1.15%  '; AsText = $false },
    @{ Ref = 'D51'; Value = '1.78'; AsText = $true },
    @{ Ref = 'E51'; Value = '  +1.68%  '; AsText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    if ($u.AsText) {
        $cell.NumberFormat = '@'
    }
    $cell.Value = $u.Value
}
